# fixed jackknife bug in dmds, new version
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dual Mode Desorption Analysis")

$data = @(
    @(8,  0.241333352, 41.48924079, 0.01),
    @(9,  0.600763584, 62.79313671, 0.01),
    @(10, 1.04806673,  77.9590348,  0.01),
    @(11, 1.466095481, 88.08019013, 0.01),
    @(12, 1.951571285, 96.61909374, 0.01),
    @(13, 2.499847618, 105.7659302, 0.01),
    @(14, 3.142683031, 114.9523482, 0.01),
    @(15, 2.60974313,  111.9833899, 0.01),
    @(16, 1.199714642, 98.93545196, 0.01),
    @(17, 0.575992209, 81.21343867, 0.01),
    @(18, 0.30991402,  66.62816092, 0.01),
    @(19, 0.145032338, 50.92125421, 0.01)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
}
